$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Overview" sheet: row 4 (fcee92e4...), column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G4").Value = "2016-10-10 09:31:35"

# "zh-cn" sheet: row 4 (fcee92e4...), column H = "Correspond Handoff Datetime", column K = "Correspond Handback DateTime"
$wsZhCn.Range("H4").Value = "2016-10-10 09:31:25"
$wsZhCn.Range("K4").Value = "2016-10-10 09:32:13"

# "de-de" sheet: row 4 (fcee92e4...), column K = "Correspond Handback DateTime"
$wsDeDe.Range("K4").Value = "2016-10-10 09:32:29"
